# (CDV) Trained rnn_010, minor edits to other files.
#
# 1. Correct the "Other Comments" text for rnn_006 (row 15): it was a
#    copy/paste of the rnn_005 comment (Tsim=1200) but rnn_006 actually
#    used Tsim=1000, matching rnn_008/rnn_009.
# 2. Log the newly trained rnn_010 run as a new row (row 23, following the
#    same "data row / blank spacer row" pattern used by the rest of the
#    table), copying the formatting from the rnn_009 row (21) so the new
#    row's styles (bold/aligned Name & Other Comments cells) match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Fix rnn_006's "Other Comments" (Tsim=1200 -> Tsim=1000) ---------
$ws.Range("I15").Value = "random uniform noise added to hidden states, excluding the first timestep, dataset normalized (featurewise) with sample mean and std, Tsim=1000"

# --- 2. Add the rnn_010 row -----------------------------------------------
# Copy formatting from the rnn_009 row (21), which has the same
# Folder/Type/Learning Rate/T_input/T_output/dt_rnn shape we need.
$ws.Range("A21:I21").Copy()
$ws.Range("A23").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("A23").Value = "colab"
$ws.Range("B23").Value = "rnn_010"
$ws.Range("C23").Value = "[64]"
$ws.Range("D23").Value = "GRU"
$ws.Range("E23").Value = "LR sigmoid (warmup 20, expected 50)"
$ws.Range("F23").Value = 90
$ws.Range("G23").Value = 90
$ws.Range("H23").Value = 0.1
$ws.Range("I23").Value = "random uniform noise (stddev=1e-3) added to hidden states, excluding the first timestep, dataset normalized (featurewise) with sample mean and std, Tsim=1000"
